$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Relabel row 4 header first: Constant -> r2
$ws.Range("A4").Value = "r2"

# Update column B (FFR Lag, LF Lag) first
$ws.Range("B2").Value = "0.68***"
$ws.Range("B3").Value = "1.246***"

# Update column C (FFR Lag, LF Lag)
$ws.Range("C2").Value = "0.382**"
$ws.Range("C3").Value = "0.837***"

# Update row 4 values: make B4/C4 numeric (was text strings)
$ws.Range("B4").Value = 0.657595502768914
$ws.Range("C4").Value = 0.5552751214566477

# Remove the old row 5 (r2_adj row) entirely
$ws.Rows("5:5").Delete()
